$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 147, shifting existing rows 147-175 down to 148-176
$ws.Rows.Item(147).Insert()

# Populate the new row 147 with the new record's data.
# Columns A,B,C,E-J are constant across this block of rows; copy them from row 148 (old row 147).
$ws.Cells.Item(147, 1).Value2 = $ws.Cells.Item(148, 1).Value2    # A
$ws.Cells.Item(147, 2).Value2 = $ws.Cells.Item(148, 2).Value2    # B
$ws.Cells.Item(147, 3).Value2 = $ws.Cells.Item(148, 3).Value2    # C

$ws.Cells.Item(147, 4).Value2 = 44543                            # D (Fecha)
$ws.Cells.Item(147, 4).NumberFormat = $ws.Cells.Item(148, 4).NumberFormat

$ws.Cells.Item(147, 5).Value2 = $ws.Cells.Item(148, 5).Value2    # E
$ws.Cells.Item(147, 6).Value2 = $ws.Cells.Item(148, 6).Value2    # F
$ws.Cells.Item(147, 7).Value2 = $ws.Cells.Item(148, 7).Value2    # G
$ws.Cells.Item(147, 8).Value2 = $ws.Cells.Item(148, 8).Value2    # H
$ws.Cells.Item(147, 9).Value2 = $ws.Cells.Item(148, 9).Value2    # I
$ws.Cells.Item(147, 10).Value2 = $ws.Cells.Item(148, 10).Value2  # J

$ws.Cells.Item(147, 11).Value2 = "Tahití"                       # K Variedad
$ws.Cells.Item(147, 12).Value2 = "Primera"                      # L Calidad
$ws.Cells.Item(147, 13).Value2 = 200                             # M Volumen
$ws.Cells.Item(147, 14).Value2 = 28000                           # N Precio minimo
$ws.Cells.Item(147, 15).Value2 = 29000                           # O Precio maximo
$ws.Cells.Item(147, 16).Value2 = 28500                           # P Precio promedio ponderado
$ws.Cells.Item(147, 17).Value2 = "$/caja 24 kilos"               # Q Unidad de comercializacion
$ws.Cells.Item(147, 18).Value2 = "Perú"                          # R Origen
$ws.Cells.Item(147, 19).Value2 = 1188                            # S Precio $/Kg
$ws.Cells.Item(147, 20).Value2 = 24                              # T Kg / unidad
